$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('B73_ISU','B73','ERS3371164','B73_ERS3371164'),
    @('B73_Ab10_ISU','B73Ab10','ERS4036196','B73Ab10_ERS4036196'),
    @('B97_ISU','B97','ERS3120542','B97_ERS3120542'),
    @('CML103_ISU','CML103','ERS3120545','CML103_ERS3120545'),
    @('CML228_ISU','CML228','ERS3120546','CML228_ERS3120546'),
    @('CML247_ISU','CML247','ERS3120547','CML247_ERS3120547'),
    @('CML277_ISU','CML277','ERS3120548','CML277_ERS3120548'),
    @('CML322_ISU','CML322','ERS3120549','CML322_ERS3120549'),
    @('CML333_ISU','CML333','ERS3120550','CML333_ERS3120550'),
    @('CML52_ISU','CML52','ERS3120543','CML52_ERS3120543'),
    @('CML69_ISU','CML69','ERS3120544','CML69_ERS3120544'),
    @('HP301_ISU','Hp301','ERS3120551','Hp301_ERS3120551'),
    @('IL14H_ISU','Il14H','ERS3120552','Il14H_ERS3120552'),
    @('Ki11_ISU','Ki11','ERS3120554','Ki11_ERS3120554'),
    @('Ki3_ISU','Ki3','ERS3120553','Ki3_ERS3120553'),
    @('Ky21_ISU','Ky21','ERS3120555','Ky21_ERS3120555'),
    @('M162W_ISU','M162W','ERS3120556','M162W_ERS3120556'),
    @('MS37W_ISU','M37W','ERS3120557','M37W_ERS3120557'),
    @('Mo18W_ISU','Mo18W','ERS3120558','Mo18W_ERS3120558'),
    @('Ms71_ISU','MS71','ERS3120559','MS71_ERS3120559'),
    @('NC350_ISU','NC350','ERS3120560','NC350_ERS3120560'),
    @('NC358_ISU','NC358','ERS3120561','NC358_ERS3120561'),
    @('OH43_ISU','Oh43','ERS3120562','Oh43_ERS3120562'),
    @('OH7B_ISU','Oh7B','ERS3120563','Oh7B_ERS3120563'),
    @('P39_ISU','P39','ERS3120564','P39_ERS3120564'),
    @('TX303_ISU','Tx303','ERS3120565','Tx303_ERS3120565'),
    @('TZi8_ISU','Tzi8','ERS3120566','Tzi8_ERS3120566'),
)

$startRow = 1500
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Match the author's final view/selection state: scrolled down so row 1500
# is the top visible row, with D1507 as the active/selected cell.
[void]$ws.Range("D1507").Select()
$excel.ActiveWindow.ScrollRow = 1500
$excel.ActiveWindow.ScrollColumn = 1
